# Update database and shift reported periods by one year (kemina income, yearly, dollar)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers (D..H) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish-date headers (D..H) ---
$ws.Range("D9").Value = "1399-04-30 (7)"
$ws.Range("E9").Value = "1400-04-28 (8)"
$ws.Range("F9").Value = "1401-04-29 (9)"
$ws.Range("G9").Value = "1402-02-29 (9)"
$ws.Range("H9").Value = "1402-02-29"

# --- Row 11: فروش ---
$ws.Range("D11").Value = 7551
$ws.Range("E11").Value = 8020
$ws.Range("F11").Value = 9882
$ws.Range("G11").Value = 12444
$ws.Range("H11").Value = 14264

# --- Row 12: بهای تمام شده کالای فروش رفته ---
$ws.Range("D12").Value = -5442
$ws.Range("E12").Value = -6040
$ws.Range("F12").Value = -5851
$ws.Range("G12").Value = -7956
$ws.Range("H12").Value = -11544

# --- Row 13: سود (زیان) ناخالص ---
$ws.Range("D13").Value = 2109
$ws.Range("E13").Value = 1980
$ws.Range("F13").Value = 4031
$ws.Range("G13").Value = 4488
$ws.Range("H13").Value = 2720

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی ---
$ws.Range("D14").Value = -612
$ws.Range("E14").Value = -541
$ws.Range("F14").Value = -420
$ws.Range("G14").Value = -650
$ws.Range("H14").Value = -970

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی) ---
$ws.Range("D15").Value = -44
$ws.Range("E15").Value = -339
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = 46
$ws.Range("E16").Value = -160
$ws.Range("F16").Value = -33
$ws.Range("G16").Value = -299
$ws.Range("H16").Value = -1074

# --- Row 17: سود (زیان) عملیاتی ---
$ws.Range("D17").Value = 1499
$ws.Range("E17").Value = 939
$ws.Range("F17").Value = 3578
$ws.Range("G17").Value = 3539
$ws.Range("H17").Value = 675

# --- Row 18: هزینه های مالی ---
$ws.Range("D18").Value = -301
$ws.Range("E18").Value = -257
$ws.Range("F18").Value = -475
$ws.Range("G18").Value = -224
$ws.Range("H18").Value = -936

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = -39
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 323

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 1198
$ws.Range("E20").Value = 643
$ws.Range("F20").Value = 3131
$ws.Range("G20").Value = 3340
$ws.Range("H20").Value = 62

# --- Row 21: مالیات ---
$ws.Range("D21").Value = -37
$ws.Range("E21").Value = -370
$ws.Range("F21").Value = -435
$ws.Range("G21").Value = -431
$ws.Range("H21").Value = -11

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 1161
$ws.Range("E22").Value = 273
$ws.Range("F22").Value = 2696
$ws.Range("G22").Value = 2909
$ws.Range("H22").Value = 51

# --- Row 24: سود (زیان) خالص ---
$ws.Range("D24").Value = 1161
$ws.Range("E24").Value = 273
$ws.Range("F24").Value = 2696
$ws.Range("G24").Value = 2909
$ws.Range("H24").Value = 51

# --- Row 26: سرمایه ---
$ws.Range("D26").Value = 2127
$ws.Range("E26").Value = 2791
$ws.Range("F26").Value = 1583
$ws.Range("G26").Value = 4825
$ws.Range("H26").Value = 3607
